$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to keep a text representation (e.g. "1.017")
    # instead of being auto-converted to a number by Excel, while
    # restoring the original cell style afterwards so no visible
    # formatting/style change is introduced.
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '27.920.86'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.861.64'
$ws.Range('E3').Value = '  -0.22%  '
Set-TextValue 'D4' '1.017'
$ws.Range('E4').Value = '  -1.85%  '
Set-TextValue 'D5' '321.52'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('E6').Value = '  -1.60%  '
Set-TextValue 'D7' '0.4340'
$ws.Range('E7').Value = '  -1.62%  '
Set-TextValue 'D8' '0.3806'
$ws.Range('E8').Value = '  +0.09%  '
Set-TextValue 'D9' '0.07441'
$ws.Range('E9').Value = '  -0.28%  '
Set-TextValue 'D10' '0.8885'
$ws.Range('E10').Value = '  +0.47%  '
Set-TextValue 'D11' '21.77'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '1.868.49'
$ws.Range('E12').Value = '  -0.27%  '
Set-TextValue 'D13' '6.789'
$ws.Range('E13').Value = '  +0.50%  '
Set-TextValue 'D14' '5.514'
$ws.Range('E14').Value = '  -0.78%  '
Set-TextValue 'D15' '0.07151'
$ws.Range('E15').Value = '  -0.85%  '
Set-TextValue 'D16' '88.50'
$ws.Range('E16').Value = '  +5.47%  '
$ws.Range('E17').Value = '  -1.67%  '
Set-TextValue 'D18' '0.000009070'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('E19').Value = '  -1.71%  '
Set-TextValue 'D20' '15.58'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '27.963.79'
$ws.Range('E21').Value = '  +0.77%  '
Set-TextValue 'D22' '5.288'
$ws.Range('E22').Value = '  -0.67%  '
Set-TextValue 'D23' '11.25'
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').Value = '2.094.75'
$ws.Range('E24').Value = '  -0.20%  '
Set-TextValue 'D25' '2.033'
$ws.Range('E25').Value = '  +4.74%  '
Set-TextValue 'D26' '157.14'
$ws.Range('E26').Value = '  -0.57%  '
Set-TextValue 'D27' '18.74'
$ws.Range('E27').Value = '  -0.62%  '
Set-TextValue 'D28' '5.439'
$ws.Range('E28').Value = '  +2.11%  '
Set-TextValue 'D29' '2.024'
$ws.Range('E29').Value = '  +1.43%  '
Set-TextValue 'D30' '120.60'
$ws.Range('E30').Value = '  +2.51%  '
Set-TextValue 'D31' '0.08997'
$ws.Range('E31').Value = '  -1.09%  '
Set-TextValue 'D32' '1.244'
$ws.Range('E32').Value = '  +2.27%  '
Set-TextValue 'D33' '0.7772'
$ws.Range('E33').Value = '  +1.02%  '
Set-TextValue 'D34' '4.606'
$ws.Range('E34').Value = '  +0.51%  '
Set-TextValue 'D35' '2.923'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D36' '1.151'
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D37' '1.018'
$ws.Range('E37').Value = '  -1.55%  '
Set-TextValue 'D38' '0.01981'
$ws.Range('E38').Value = '  -0.54%  '
Set-TextValue 'D39' '0.05335'
$ws.Range('E39').Value = '  -0.36%  '
Set-TextValue 'D40' '2.887'
$ws.Range('E40').Value = '  +1.82%  '
Set-TextValue 'D41' '0.5226'
$ws.Range('E41').Value = '  +0.48%  '
Set-TextValue 'D42' '7.035'
$ws.Range('E42').Value = '  +2.82%  '
Set-TextValue 'D43' '0.1687'
$ws.Range('E43').Value = '  -0.40%  '
Set-TextValue 'D44' '8.813'
$ws.Range('E44').Value = '  +1.24%  '
Set-TextValue 'D45' '111.23'
$ws.Range('E45').Value = '  +1.57%  '
Set-TextValue 'D46' '10.75'
$ws.Range('E46').Value = '  +1.06%  '
Set-TextValue 'D47' '0.4774'
$ws.Range('E47').Value = '  +1.62%  '
Set-TextValue 'D48' '1.719'
$ws.Range('E48').Value = '  -0.70%  '
Set-TextValue 'D49' '0.06495'
$ws.Range('E49').Value = '  +1.01%  '
Set-TextValue 'D50' '1.018'
$ws.Range('E50').Value = '  -1.66%  '
Set-TextValue 'D51' '1.888'
$ws.Range('E51').Value = '  +0.45%  '
